$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a fresh blank row at position 2. This pushes the existing header
#    (row 1) content down to... no: inserting AT row 2 keeps row 1 in place
#    and shifts old rows 2..15 down to 3..16 (preserving their per-row
#    formatting, e.g. old row 8 -> row 9 keeps its custom row height/style).
# ---------------------------------------------------------------------------
$ws.Rows("2").Insert()

# ---------------------------------------------------------------------------
# 2. The old header row (currently still row 1: Algorithms/Notes/TC-Average/
#    TC-WC/SC-Average.../SC-WC) moves logically to row 2. Populate row 2 with
#    that header data, but with the new "Time for each node" / blank / O(height)... / blank
#    values for the TC/SC columns, then clear out row 1 and refill it with the
#    new short sub-header row (TC-Average / TC-WC / SC-AC / SC-WC).
#    New-string cells are written in SC-AC -> Time for each node -> O(height)...
#    order so sharedStrings.xml gets the same append order as the target file.
# ---------------------------------------------------------------------------

# -- Row 1 (new short sub-header) --
$ws.Range("A1").Clear()
$ws.Range("B1").Clear()
$ws.Range("C1").Value = "TC-Average"
$ws.Range("D1").Value = "TC-WC"
$ws.Range("E1").Value = "SC-AC"
$ws.Range("F1").Value = "SC-WC"

$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").VerticalAlignment = -4160
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").VerticalAlignment = -4160

$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").VerticalAlignment = -4160
$ws.Range("E1").WrapText = $true

$ws.Rows("1").RowHeight = 17

# -- Row 2 (was row 1) --
$ws.Range("A2").Value = "Algorithms"
$ws.Range("B2").Value = "Notes"
$ws.Range("C2").Value = "Time for each node"
$ws.Range("E2").Value = "O(height), 粉色path的深度"

$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").VerticalAlignment = -4160

$ws.Range("B2").Font.Bold = $true

$ws.Range("C2:F2").HorizontalAlignment = -4108
$ws.Range("C2:D2").Merge()
$ws.Range("E2:F2").Merge()

# ---------------------------------------------------------------------------
# 3. Append the new last row: "Binary Tree - isSymmetric"
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = "Binary Tree - isSymmetric"

# ---------------------------------------------------------------------------
# 4. Column width tweaks (D, E, F)
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 9.33203125
$ws.Columns("E").ColumnWidth = 10.83203125
$ws.Columns("F").ColumnWidth = 14

# ---------------------------------------------------------------------------
# 5. Freeze panes: freeze first column + first two rows, top-left visible
#    cell of the scrolling area is B3.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("B3").Select()
$win.FreezePanes = $true

# ---------------------------------------------------------------------------
# 6. Final selection
# ---------------------------------------------------------------------------
$ws.Range("C17").Select()
